# Add PAGE input perturbation columns for "2080" shock-year to both the
# CH4 table (rows 2-13) and the N2O table (rows 16-27), mirroring the
# existing 2010..2060 blocks which each occupy 6 columns (Year + 5 model
# columns: IMAGE, MESSAGE, MiniCAM, MERGE, Policy) starting at A, H, O, V,
# AC, AJ, AQ. The new "2080" block goes in the next free slot: AX:BC.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Clone formatting (fills/alignment) for the new block's header row
#    and data rows from the existing "2060" block (AQ:AV), for both
#    tables. This brings over styles s="4" (header) and s="5" (data)
#    without having to hard-code style indices.
# ---------------------------------------------------------------------
$ws.Range("AQ3:AV3").Copy($ws.Range("AX3:BC3"))
$ws.Range("AQ4:AV13").Copy($ws.Range("AX4:BC13"))

$ws.Range("AQ17:AV17").Copy($ws.Range("AX17:BC17"))
$ws.Range("AQ18:AV27").Copy($ws.Range("AX18:BC27"))

# ---------------------------------------------------------------------
# 2. New section labels (row 2 / row 16) naming the new shock-year block.
# ---------------------------------------------------------------------
$ws.Range("AX2").Value = "PAGE input: CH4 Shock, 2080"
$ws.Range("AX16").Value = "PAGE input: N2O Shock, 2080"

# ---------------------------------------------------------------------
# 3. Fill in the actual Year + data values for each table.
#    Columns: AX=Year, AY=IMAGE, AZ=MESSAGE, BA=MiniCAM, BB=MERGE, BC=Policy
# ---------------------------------------------------------------------

$ch4 = @(
    @(4,  "2010", "0","0","0","0","0"),
    @(5,  "2020", "0","0","0","0","0"),
    @(6,  "2030", "0","0","0","0","0"),
    @(7,  "2040", "0","0","0","0","0"),
    @(8,  "2050", "0","0","0","0","0"),
    @(9,  "2060", "0","0","0","0","0"),
    @(10, "2080", "7.2240264871450524E-5", "7.5769316062878004E-5", "7.8418998831559543E-5", "7.6237339676625473E-5", "9.3649217227735539E-5"),
    @(11, "2100", "3.1299133349849127E-6",  "3.4013421684231383E-6",  "3.5168716620059827E-6",  "3.2861195249855866E-6",  "4.2130350609481668E-6"),
    @(12, "2200", "5.1471513939915784E-10", "5.6612308174663899E-10", "5.8461658847264171E-10", "5.3940619304526645E-10", "7.00520939034277E-10"),
    @(13, "2300", "7.1409544943890069E-13", "7.8492767840998567E-13", "8.1090689718621434E-13", "7.4784622938750545E-13", "9.7144514654701197E-13")
)

$n2o = @(
    @(18, "2010", "0","0","0","0","0"),
    @(19, "2020", "0","0","0","0","0"),
    @(20, "2030", "0","0","0","0","0"),
    @(21, "2040", "0","0","0","0","0"),
    @(22, "2050", "0","0","0","0","0"),
    @(23, "2060", "0","0","0","0","0"),
    @(24, "2080", "3.6132447295342409E-4", "3.415610660511957E-4",  "3.371325169616462E-4",  "3.4770880718886757E-4", "3.5550284434211548E-4"),
    @(25, "2100", "2.1598616098579383E-4", "2.0637577727042057E-4", "1.9733902078532316E-4", "2.0583216002156179E-4", "2.1393167721237949E-4"),
    @(26, "2200", "8.752898405308795E-5",  "8.4698988014204611E-5", "7.8355453610929883E-5", "8.2606011835180122E-5", "8.7322830481476044E-5"),
    @(27, "2300", "5.3976699967039909E-5", "5.241843039205385E-5",  "4.8054081898007972E-5", "5.0805233537953853E-5", "5.3958351930005755E-5")
)

foreach ($row in $ch4) {
    $r = $row[0]
    $ws.Cells.Item($r, 50).Value = [double]$row[1]
    $ws.Cells.Item($r, 51).Value = [double]$row[2]
    $ws.Cells.Item($r, 52).Value = [double]$row[3]
    $ws.Cells.Item($r, 53).Value = [double]$row[4]
    $ws.Cells.Item($r, 54).Value = [double]$row[5]
    $ws.Cells.Item($r, 55).Value = [double]$row[6]
}

foreach ($row in $n2o) {
    $r = $row[0]
    $ws.Cells.Item($r, 50).Value = [double]$row[1]
    $ws.Cells.Item($r, 51).Value = [double]$row[2]
    $ws.Cells.Item($r, 52).Value = [double]$row[3]
    $ws.Cells.Item($r, 53).Value = [double]$row[4]
    $ws.Cells.Item($r, 54).Value = [double]$row[5]
    $ws.Cells.Item($r, 55).Value = [double]$row[6]
}

# ---------------------------------------------------------------------
# 4. Mirror the final selection/scroll position from the edit.
# ---------------------------------------------------------------------
$ws.Range("AZ14").Select()
